$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1488.909
$ws.Range("I70").Value = 1496.4445
$ws.Range("J70").Value = 1483.6923
$ws.Range("K70").Value = 4489.333500000001
$ws.Range("L70").Value = 4451.0769
$ws.Range("M70").Value = -4219.333500000001
$ws.Range("N70").Value = -4991.0769
$ws.Range("H73").Value = 1488.909
$ws.Range("I73").Value = 1496.4445
$ws.Range("J73").Value = 1483.6923
$ws.Range("K73").Value = 4489.333500000001
$ws.Range("L73").Value = 4451.0769
$ws.Range("M73").Value = -3553.333500000001
$ws.Range("N73").Value = -6323.0769

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3711.2104
$ws.Range("I132").Value = 2284.875
$ws.Range("J132").Value = 4748.5454
$ws.Range("K132").Value = 6854.625
$ws.Range("L132").Value = 14245.6362
$ws.Range("M132").Value = -4324.625
$ws.Range("N132").Value = -19305.6362
$ws.Range("H134").Value = 48763.625
$ws.Range("J134").Value = 48763.625
$ws.Range("L134").Value = 48763.625
$ws.Range("N134").Value = -58903.625

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 200
$ws.Range("I22").Value = 200
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 200
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -27
$ws.Range("H86").Value = 1659.8077
$ws.Range("I86").Value = 1734.975
$ws.Range("J86").Value = 1409.25
$ws.Range("K86").Value = 1734.975
$ws.Range("L86").Value = 1409.25
$ws.Range("M86").Value = -611.9749999999999
$ws.Range("N86").Value = -3655.25
$ws.Range("H89").Value = 1659.8077
$ws.Range("I89").Value = 1734.975
$ws.Range("J89").Value = 1409.25
$ws.Range("K89").Value = 8674.875
$ws.Range("L89").Value = 7046.25
$ws.Range("M89").Value = -3058.875
$ws.Range("N89").Value = -18278.25
$ws.Range("H118").Value = 57388.11
$ws.Range("J118").Value = 57388.11
$ws.Range("L118").Value = 57388.11
$ws.Range("N118").Value = -60702.11
$ws.Range("H134").Value = 5453.7393
$ws.Range("I134").Value = 4396.4443
$ws.Range("K134").Value = 13189.3329
$ws.Range("M134").Value = -10654.3329

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 2949.5334
$ws.Range("I80").Value = 3100
$ws.Range("J80").Value = 2911.9167
$ws.Range("K80").Value = 9300
$ws.Range("L80").Value = 8735.750100000001
$ws.Range("M80").Value = -8364
$ws.Range("N80").Value = -10607.7501
$ws.Range("H83").Value = 2949.5334
$ws.Range("I83").Value = 3100
$ws.Range("J83").Value = 2911.9167
$ws.Range("K83").Value = 27900
$ws.Range("L83").Value = 26207.2503
$ws.Range("M83").Value = -23220
$ws.Range("N83").Value = -35567.2503
$ws.Range("H92").Value = 1367.6666
$ws.Range("I92").Value = 1300
$ws.Range("J92").Value = 1401.5
$ws.Range("K92").Value = 3900
$ws.Range("L92").Value = 4204.5
$ws.Range("M92").Value = -2652
$ws.Range("N92").Value = -6700.5
$ws.Range("H121").Value = 223.25
$ws.Range("J121").Value = 216.5
$ws.Range("L121").Value = 649.5
$ws.Range("N121").Value = -3269.5

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 40000
$ws.Range("J32").Value = 40000
$ws.Range("L32").Value = 40000
$ws.Range("N32").Value = -40592
$ws.Range("H42").Value = 60120
$ws.Range("J42").Value = 60120
$ws.Range("L42").Value = 60120
$ws.Range("N42").Value = -61090
$ws.Range("H80").Value = 7283.769
$ws.Range("I80").Value = 26149.5
$ws.Range("K80").Value = 26149.5
$ws.Range("M80").Value = -25151.5
$ws.Range("H83").Value = 7283.769
$ws.Range("I83").Value = 26149.5
$ws.Range("K83").Value = 130747.5
$ws.Range("M83").Value = -125755.5
$ws.Range("H115").Value = 60120
$ws.Range("J115").Value = 60120
$ws.Range("L115").Value = 60120
$ws.Range("N115").Value = -62470

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 16611.354
$ws.Range("I61").Value = 30530.285
$ws.Range("J61").Value = 6868.1
$ws.Range("K61").Value = 30530.285
$ws.Range("L61").Value = 6868.1
$ws.Range("M61").Value = -30328.285
$ws.Range("N61").Value = -7272.1
$ws.Range("H68").Value = 1375
$ws.Range("I68").Value = 1375
$ws.Range("K68").Value = 1375
$ws.Range("M68").Value = -626
$ws.Range("H71").Value = 1375
$ws.Range("I71").Value = 1375
$ws.Range("K71").Value = 6875
$ws.Range("M71").Value = -3131
$ws.Range("H113").Value = 16611.354
$ws.Range("I113").Value = 30530.285
$ws.Range("J113").Value = 6868.1
$ws.Range("K113").Value = 30530.285
$ws.Range("L113").Value = 6868.1
$ws.Range("M113").Value = -28360.285
$ws.Range("N113").Value = -11208.1

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 51000
$ws.Range("J27").Value = 51000
$ws.Range("L27").Value = 51000
$ws.Range("N27").Value = -51138
$ws.Range("H57").Value = 39333.332
$ws.Range("J57").Value = 39333.332
$ws.Range("L57").Value = 39333.332
$ws.Range("N57").Value = -40841.332
$ws.Range("H62").Value = 4071.647
$ws.Range("J62").Value = 4135.3335
$ws.Range("L62").Value = 4135.3335
$ws.Range("N62").Value = -5383.3335
$ws.Range("H65").Value = 4071.647
$ws.Range("J65").Value = 4135.3335
$ws.Range("L65").Value = 20676.6675
$ws.Range("N65").Value = -26916.6675
$ws.Range("H81").Value = 1493.75
$ws.Range("I81").Value = 991.6667
$ws.Range("K81").Value = 1983.3334
$ws.Range("M81").Value = -922.3334
$ws.Range("H84").Value = 1493.75
$ws.Range("I84").Value = 991.6667
$ws.Range("K84").Value = 9916.666999999999
$ws.Range("M84").Value = -4612.666999999999
$ws.Range("H86").Value = 23800
$ws.Range("J86").Value = 23800
$ws.Range("L86").Value = 23800
$ws.Range("N86").Value = -26046
$ws.Range("H89").Value = 23800
$ws.Range("J89").Value = 23800
$ws.Range("L89").Value = 119000
$ws.Range("N89").Value = -130232
$ws.Range("I136").Value = 2426.56
$ws.Range("J136").Value = 7657.2915
$ws.Range("K136").Value = 7279.68
$ws.Range("L136").Value = 22971.8745
$ws.Range("M136").Value = -4729.68
$ws.Range("N136").Value = -28071.8745
